# Hortaliza, Macroferia Regional de Talca - Arveja Verde
# A new daily price record is inserted at row 59, pushing the existing
# rows 59-126 down to 60-127 (dimension grows from A1:R126 to A1:R127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 59, shifting rows 59:126 -> 60:127.
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new observation.
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "Macroferia Regional de Talca"
$ws.Range("C59").Value = "Maule"
$ws.Range("D59").Value = 44902
$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 100112022
$ws.Range("G59").Value = "Arveja Verde"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 300
$ws.Range("K59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("M59").Value = 20000
$ws.Range("N59").Value = "`$/saco 25 kilos"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 800
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
